# TemplateDiem.xlsx fix: lop hoc + dot thi + import diem thi
# Replace "ngaythi / giothi / phongthi" columns (I/J/K) with the
# laptrinhc/msword/msexcel/mspowerpoint block that used to live in L:O,
# shifting it three columns to the left (I:L) and dropping the old M:O data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Cells.Item(1, 9).Value  = "laptrinhc"     # I1 (was ngaythi)
$ws.Cells.Item(1, 10).Value = "msword"        # J1 (was giothi)
$ws.Cells.Item(1, 11).Value = "msexcel"       # K1 (was phongthi)
$ws.Cells.Item(1, 12).Value = "mspowerpoint"  # L1 (unchanged content, but now re-set)
$ws.Range("M1:O1").Clear()

# ---- Row 2 (Nguyen Van A) ----
$ws.Cells.Item(2, 9).ClearFormats()
$ws.Cells.Item(2, 9).Value  = 10   # I2 (was date 44167 w/ style)
$ws.Cells.Item(2, 10).Value = 10   # J2 (was text "18h")
$ws.Cells.Item(2, 11).Value = 10   # K2 (was 16.5)
$ws.Cells.Item(2, 12).Value = 10   # L2 (already 10)
$ws.Range("M2:O2").Clear()

# ---- Row 3 (Nguyen Van B) ----
$ws.Cells.Item(3, 9).ClearFormats()
$ws.Cells.Item(3, 9).Value  = 0    # I3 (was date 1206370 w/ style)
$ws.Cells.Item(3, 10).Value = 0    # J3 (was text "12h")
$ws.Cells.Item(3, 11).Value = 0    # K3 (was 162)
$ws.Cells.Item(3, 12).Value = 0    # L3 (already 0)
$ws.Range("M3:O3").Clear()

# ---- View state: scroll right a bit and move the active selection ----
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("K5").Select()
